# Applies the crypto-price / volume refresh for Mon Nov 18 01:09:31 UTC 2024 run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "90.046.78"
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("D3").Value = "3.099.61"
$ws.Range("E3").Value = "  -1.51%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.47"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +8.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "622.50"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.07"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -4.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.364"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.72%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "3.096.27"
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.711"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -5.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.202"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000246"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.81"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.40"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.47%  "
$ws.Range("D16").Value = "89.759.35"
$ws.Range("E16").Value = "  -0.89%  "
$ws.Range("D17").Value = "3.678.62"
$ws.Range("E17").Value = "  -1.82%  "
$ws.Range("D18").Value = "3.089.14"
$ws.Range("E18").Value = "  -2.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.78"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("E20").Value = "  -0.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.02"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -5.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "438.91"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -8.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.55"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +6.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.90"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -4.47%  "
$ws.Range("E25").Value = "  -6.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "87.52"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -8.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.86"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.37%  "
$ws.Range("D28").Value = "3.241.24"
$ws.Range("E28").Value = "  -4.50%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("B30").Value = "Cronos"
$ws.Range("C30").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.165"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.47%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.99"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -5.11%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.202"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -6.80%  "
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.947"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +24.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.79"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -9.00%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.150"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.03%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.20"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.52%  "
$ws.Range("B37").Value = "MantraDAO"
$ws.Range("C37").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.01"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +31.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "492.99"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.60"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.89"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.98%  "
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0922"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +6.40%  "
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.27"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.10"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.66%  "
$ws.Range("B44").Value = "PolygonEcosystemToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.405"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -7.00%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "156.33"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.35%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.693"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.36%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.87"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -5.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.63"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.70%  "
$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.38"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -8.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.31"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.66%  "
